# Edit script for SPER.xlsx
# 1) Update raw "Temps Passé Période (h)" entries on the Octobre sheet
# 2) Flip the matching "Terminé" (Oui/Non) flags on the Octobre sheet
# 3) Swap the "Temps Passé Cumulé" / "Ecart de Charge Total" columns (E <-> F)
#    on the Synthèse sheet, and simplify the Avancement formula
# 4) Leave the final selection on Synthèse!I4

$wb = $excel.ActiveWorkbook

# ---- 1) & 2) Octobre worksheet raw-data edits ----
$octobre = $wb.Worksheets.Item("Octobre")

$octobre.Range("E2").Value = 1
$octobre.Range("E3").Value = 3
$octobre.Range("E5").Value = 2
$octobre.Range("E6").Value = 8
$octobre.Range("E7").Value = 3
$octobre.Range("E8").Value = 3

$octobre.Range("N3").Value = "Oui"
$octobre.Range("N5").Value = "Oui"
$octobre.Range("N6").Value = "Oui"
$octobre.Range("N7").Value = "Oui"
$octobre.Range("N8").Value = "Oui"

# ---- 3) Synthèse worksheet: swap columns E and F ----
$synthese = $wb.Worksheets.Item("Synthèse")

# Headers (row 3)
$eHeader = $synthese.Cells.Item(3, 5).Value2
$fHeader = $synthese.Cells.Item(3, 6).Value2
$synthese.Cells.Item(3, 5).Value2 = $fHeader
$synthese.Cells.Item(3, 6).Value2 = $eHeader

# Formulas (rows 4-9)
for ($r = 4; $r -le 9; $r++) {
    $eFormula = $synthese.Cells.Item($r, 5).Formula
    $fFormula = $synthese.Cells.Item($r, 6).Formula
    $synthese.Cells.Item($r, 5).Formula = $fFormula
    $synthese.Cells.Item($r, 6).Formula = $eFormula
}

# Avancement column formula no longer needs E, just F
$synthese.Range("I4").Formula = "=IFERROR(F4/D4,0)"
$synthese.Range("I5:I9").Formula = "=IFERROR(F5/D5,0)"

# ---- 4) Refresh charts so their cached values follow the new data ----
for ($i = 1; $i -le $synthese.ChartObjects().Count; $i++) {
    $synthese.ChartObjects($i).Chart.Refresh()
}

# ---- 5) Final selection ----
$synthese.Activate()
$synthese.Range("I4").Select()
